# Apply the edit described by the diff:
# 1. Negate the "long" (longitude) values in column B, rows 2-11.
# 2. Change the selection on Sheet1 from D1:D1048576 to D15 (single cell D15).
# 3. (Window position xWindow/yWindow change is an application-level view
#    setting not tied to worksheet content; left as-is by default behavior.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Negate column B values for rows 2 through 11
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = -1 * $cell.Value2
}

# Update the active selection to D15
$ws.Activate()
$ws.Range("D15").Select()
